$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DANH SÁCH NỢ")
$ws.Columns("L").Insert()
$ws.Columns("L").ColumnWidth = 13.28515625

$ws.Range("L1").Value = "Còn lại"
